$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.056.40"
$ws.Range("E2").Value = "'  +0.56%  "
$ws.Range("D3").Value = "'2.305.18"
$ws.Range("E3").Value = "'  +0.64%  "
$ws.Range("E4").Value = "'  -0.12%  "
$ws.Range("D5").Value = "'305.00"
$ws.Range("E5").Value = "'  +1.76%  "
$ws.Range("D6").Value = "'97.32"
$ws.Range("E6").Value = "'  +0.75%  "
$ws.Range("E7").Value = "'  -1.19%  "
$ws.Range("E8").Value = "'  -0.05%  "
$ws.Range("D9").Value = "'0.503"
$ws.Range("E9").Value = "'  -0.21%  "
$ws.Range("D10").Value = "'35.46"
$ws.Range("E10").Value = "'  -0.58%  "
$ws.Range("D11").Value = "'0.0788"
$ws.Range("E11").Value = "'  +0.25%  "
$ws.Range("D12").Value = "'18.61"
$ws.Range("E12").Value = "'  +5.09%  "
$ws.Range("E13").Value = "'  +1.68%  "
$ws.Range("E14").Value = "'  +2.46%  "
$ws.Range("D15").Value = "'2.664.14"
$ws.Range("E15").Value = "'  +0.63%  "
$ws.Range("D16").Value = "'2.311.06"
$ws.Range("E16").Value = "'  +0.94%  "
$ws.Range("E17").Value = "'  +1.16%  "
$ws.Range("D18").Value = "'42.937.95"
$ws.Range("E18").Value = "'  +0.45%  "
$ws.Range("D19").Value = "'12.59"
$ws.Range("E19").Value = "'  -0.66%  "
$ws.Range("D20").Value = "'0.0₃0900"
$ws.Range("E20").Value = "'  -0.47%  "
$ws.Range("D21").Value = "'6.05"
$ws.Range("E21").Value = "'  -0.10%  "
$ws.Range("D22").Value = "'67.58"
$ws.Range("E22").Value = "'  -0.30%  "
$ws.Range("D23").Value = "'237.05"
$ws.Range("E23").Value = "'  -1.71%  "
$ws.Range("E24").Value = "'  +2.24%  "
$ws.Range("E25").Value = "'  +0.12%  "
$ws.Range("E26").Value = "'  +0.43%  "
$ws.Range("D27").Value = "'25.06"
$ws.Range("E27").Value = "'  -0.25%  "
$ws.Range("D28").Value = "'2.19"
$ws.Range("E28").Value = "'  +7.61%  "
$ws.Range("D29").Value = "'166.17"
$ws.Range("E29").Value = "'  +0.17%  "
$ws.Range("E30").Value = "'  +0.26%  "
$ws.Range("E31").Value = "'  +0.29%  "
$ws.Range("E32").Value = "'  +0.04%  "
$ws.Range("D33").Value = "'18.22"
$ws.Range("E33").Value = "'  +6.67%  "
$ws.Range("E34").Value = "'  -0.20%  "
$ws.Range("E35").Value = "'  -7.75%  "
$ws.Range("E36").Value = "'  -0.92%  "
$ws.Range("D37").Value = "'0.0691"
$ws.Range("E37").Value = "'  +1.00%  "
$ws.Range("E38").Value = "'  +0.13%  "
$ws.Range("E39").Value = "'  +0.42%  "
$ws.Range("D40").Value = "'2.75"
$ws.Range("E40").Value = "'  +0.75%  "
$ws.Range("E41").Value = "'  -0.44%  "
$ws.Range("D42").Value = "'1.997.19"
$ws.Range("E42").Value = "'  -0.58%  "
$ws.Range("E43").Value = "'  +3.52%  "
$ws.Range("E44").Value = "'  -0.48%  "
$ws.Range("B45").Value = "'EnergySwap"
$ws.Range("C45").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'17.99"
$ws.Range("E45").Value = "'  +5.34%  "
$ws.Range("B46").Value = "'ApeXProtocol"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'2.11"
$ws.Range("E46").Value = "'  +2.21%  "
$ws.Range("E47").Value = "'  +0.49%  "
$ws.Range("D48").Value = "'53.59"
$ws.Range("E48").Value = "'  +0.98%  "
$ws.Range("D49").Value = "'2.532.13"
$ws.Range("E50").Value = "'  -7.53%  "
$ws.Range("D51").Value = "'71.82"
$ws.Range("E51").Value = "'  -0.23%  "
